# This script fixes a metal-name translation bug: the "Mercure" (Mercury)
# and "Plomb" (Lead) rows had their data rows correctly labeled in French,
# but when translated to English, the labels were swapped relative to the
# actual underlying data. The fix renames "Mercure" -> "Lead" and
# "Plomb" -> "Mercury", while keeping the originally-correct data attached
# to the correct (now-renamed) label - i.e. the numeric rows for the old
# "Mercure" row move down to the old "Plomb" row (now "Mercury"), and vice
# versa.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet1: small lookup/summary table (PARAMETRE_LIBELLE, Gironde, Loire, Seine)
# Row 3 = Mercure, Row 4 = Plomb. Swap the B/C/D values between the two rows
# and rename the labels.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$row3 = @($ws1.Range("B3").Value(), $ws1.Range("C3").Value(), $ws1.Range("D3").Value())
$row4 = @($ws1.Range("B4").Value(), $ws1.Range("C4").Value(), $ws1.Range("D4").Value())

$ws1.Range("A3").Value = "Lead"
$ws1.Range("B3").Value = $row4[0]
$ws1.Range("C3").Value = $row4[1]
$ws1.Range("D3").Value = $row4[2]

$ws1.Range("A4").Value = "Mercury"
$ws1.Range("B4").Value = $row3[0]
$ws1.Range("C4").Value = $row3[1]
$ws1.Range("D4").Value = $row3[2]

# ---------------------------------------------------------------------
# Sheet2: statistics table (ESTUARY, PARAMETRE_LIBELLE, rho, p.value,
# short_last_trend, EC_MPC_µg_gww, pvalue, median_1, median_2,
# long_term_trend, status)
# Pairs of rows (Mercure, Plomb) per estuary: (3,4) Gironde, (6,7) Loire,
# (9,10) Seine. For each pair, rename labels and swap the C..J data
# columns between the two rows (K/status column stays put).
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$pairs = @(
    @(3, 4),
    @(6, 7),
    @(9, 10)
)

foreach ($pair in $pairs) {
    $rMercure = $pair[0]
    $rPlomb = $pair[1]

    # Columns (1-based): 3=rho, 4=p.value, 5=short_last_trend, 6=EC_MPC_µg_gww,
    # 7=pvalue (numeric-looking TEXT in the source file), 8=median_1,
    # 9=median_2, 10=long_term_trend. Column 7 ("pvalue") is stored as text
    # (e.g. "0.0523") rather than a real number in the original workbook, so
    # it needs special handling on write-back to avoid Excel's automatic
    # text->number coercion.
    $dataMercure = @(
        $ws2.Cells.Item($rMercure, 3).Value(),
        $ws2.Cells.Item($rMercure, 4).Value(),
        $ws2.Cells.Item($rMercure, 5).Value(),
        $ws2.Cells.Item($rMercure, 6).Value(),
        $ws2.Cells.Item($rMercure, 7).Text,
        $ws2.Cells.Item($rMercure, 8).Value(),
        $ws2.Cells.Item($rMercure, 9).Value(),
        $ws2.Cells.Item($rMercure, 10).Value()
    )

    $dataPlomb = @(
        $ws2.Cells.Item($rPlomb, 3).Value(),
        $ws2.Cells.Item($rPlomb, 4).Value(),
        $ws2.Cells.Item($rPlomb, 5).Value(),
        $ws2.Cells.Item($rPlomb, 6).Value(),
        $ws2.Cells.Item($rPlomb, 7).Text,
        $ws2.Cells.Item($rPlomb, 8).Value(),
        $ws2.Cells.Item($rPlomb, 9).Value(),
        $ws2.Cells.Item($rPlomb, 10).Value()
    )

    # Former "Mercure" row becomes "Lead" and receives the former "Plomb" data
    $ws2.Cells.Item($rMercure, 2).Value = "Lead"
    $ws2.Cells.Item($rMercure, 3).Value = $dataPlomb[0]
    $ws2.Cells.Item($rMercure, 4).Value = $dataPlomb[1]
    $ws2.Cells.Item($rMercure, 5).Value = $dataPlomb[2]
    $ws2.Cells.Item($rMercure, 6).Value = $dataPlomb[3]
    $gCell = $ws2.Cells.Item($rMercure, 7)
    $gCell.NumberFormat = "@"
    $gCell.Value = $dataPlomb[4]
    $gCell.Style = "Normal"
    $ws2.Cells.Item($rMercure, 8).Value = $dataPlomb[5]
    $ws2.Cells.Item($rMercure, 9).Value = $dataPlomb[6]
    $ws2.Cells.Item($rMercure, 10).Value = $dataPlomb[7]

    # Former "Plomb" row becomes "Mercury" and receives the former "Mercure" data
    $ws2.Cells.Item($rPlomb, 2).Value = "Mercury"
    $ws2.Cells.Item($rPlomb, 3).Value = $dataMercure[0]
    $ws2.Cells.Item($rPlomb, 4).Value = $dataMercure[1]
    $ws2.Cells.Item($rPlomb, 5).Value = $dataMercure[2]
    $ws2.Cells.Item($rPlomb, 6).Value = $dataMercure[3]
    $gCell2 = $ws2.Cells.Item($rPlomb, 7)
    $gCell2.NumberFormat = "@"
    $gCell2.Value = $dataMercure[4]
    $gCell2.Style = "Normal"
    $ws2.Cells.Item($rPlomb, 8).Value = $dataMercure[5]
    $ws2.Cells.Item($rPlomb, 9).Value = $dataMercure[6]
    $ws2.Cells.Item($rPlomb, 10).Value = $dataMercure[7]
}
